$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text number format on numeric-looking columns for rows 35-48
$ws.Range("F35:F48").NumberFormat = "@"
$ws.Range("I35:I48").NumberFormat = "@"
$ws.Range("J35:J48").NumberFormat = "@"
$ws.Range("K35:K48").NumberFormat = "@"

# Row 35
$ws.Range("A35").Value = 'Raxmatov Ogabek'
$ws.Range("B35").Value = 'Yurisprudensiya'
$ws.Range("C35").Value = 'O''zbek tili'
$ws.Range("D35").Value = 'Kunduzgi'
$ws.Range("E35").Value = 'AD4419590'
$ws.Range("F35").Value = '50908076150018'
$ws.Range("G35").Value = 'Toshkent shahri'
$ws.Range("H35").Value = 'Yashnaobod tumani'
$ws.Range("I35").Value = '998336373784'
$ws.Range("J35").Value = '+998931984133'
$ws.Range("K35").Value = '2025-05-12'

# Row 36
$ws.Range("A36").Value = 'Boboqulova Bahora Sherzod qizi'
$ws.Range("B36").Value = 'Yurisprudensiya'
$ws.Range("C36").Value = 'O''zbek tili'
$ws.Range("D36").Value = 'Kunduzgi'
$ws.Range("E36").Value = 'AD6081128'
$ws.Range("F36").Value = '60803066050143'
$ws.Range("G36").Value = 'Samarqand viloyati'
$ws.Range("H36").Value = 'Ishtixon tumani'
$ws.Range("I36").Value = '998944292086'
$ws.Range("J36").Value = '+998944292086'
$ws.Range("K36").Value = '2025-05-12'

# Row 37
$ws.Range("A37").Value = 'Malikova Shoxidaxon Zakirjanovna'
$ws.Range("B37").Value = 'Yurisprudensiya'
$ws.Range("C37").Value = 'O''zbek tili'
$ws.Range("D37").Value = 'Kunduzgi'
$ws.Range("E37").Value = 'AD7703787'
$ws.Range("F37").Value = '41901891230031'
$ws.Range("G37").Value = 'Andijon viloyati'
$ws.Range("H37").Value = 'Andijon shahri'
$ws.Range("I37").Value = '998999767752'
$ws.Range("J37").Value = '+998999767752'
$ws.Range("K37").Value = '2025-05-13'

# Row 38
$ws.Range("A38").Value = 'Qudratov Sarvar Raximjon ogli'
$ws.Range("B38").Value = 'Yurisprudensiya'
$ws.Range("C38").Value = 'O''zbek tili'
$ws.Range("D38").Value = 'Kunduzgi'
$ws.Range("E38").Value = 'AD3368539'
$ws.Range("F38").Value = '50604075670019'
$ws.Range("G38").Value = 'Toshkent shahri'
$ws.Range("H38").Value = 'Olmazor tumani'
$ws.Range("I38").Value = '998940473807'
$ws.Range("J38").Value = '+998330073407'
$ws.Range("K38").Value = '2025-05-13'

# Row 39
$ws.Range("A39").Value = 'O''roqov Xushnudbek Xurshid o''g''li'
$ws.Range("B39").Value = 'Yurisprudensiya'
$ws.Range("C39").Value = 'O''zbek tili'
$ws.Range("D39").Value = 'Kunduzgi'
$ws.Range("E39").Value = 'AE2454790'
$ws.Range("F39").Value = '51804076540053'
$ws.Range("G39").Value = 'Toshkent shahri'
$ws.Range("H39").Value = 'Shayxontohur tumani'
$ws.Range("I39").Value = '998999840771'
$ws.Range("J39").Value = '+998999840771'
$ws.Range("K39").Value = '2025-05-14'

# Row 40
$ws.Range("A40").Value = 'Iskandarov Sarvarbek Iskandar o''g''li'
$ws.Range("B40").Value = 'Yurisprudensiya'
$ws.Range("C40").Value = 'O''zbek tili'
$ws.Range("D40").Value = 'Kunduzgi'
$ws.Range("E40").Value = 'AD4799601'
$ws.Range("F40").Value = '52911076170021'
$ws.Range("G40").Value = 'Samarqand viloyati'
$ws.Range("H40").Value = 'Kattaqoʻrgʻon tumani'
$ws.Range("I40").Value = '998944475702'
$ws.Range("J40").Value = '+998944475702'
$ws.Range("K40").Value = '2025-05-14'

# Row 41
$ws.Range("A41").Value = 'Numonjonov Zuhriddin Abdulhafiz ogli'
$ws.Range("B41").Value = 'Yurisprudensiya'
$ws.Range("C41").Value = 'O''zbek tili'
$ws.Range("D41").Value = 'Kunduzgi'
$ws.Range("E41").Value = 'AD7653210'
$ws.Range("F41").Value = '50612077080048'
$ws.Range("G41").Value = 'Fargona viloyati'
$ws.Range("H41").Value = 'Margʻilon tumani'
$ws.Range("I41").Value = '998911245111'
$ws.Range("J41").Value = '+998331590053'
$ws.Range("K41").Value = '2025-05-14'

# Row 42
$ws.Range("A42").Value = 'Muhammadaliyev Jahongir Mahmudjon ogli'
$ws.Range("B42").Value = 'Psixologiya'
$ws.Range("C42").Value = 'O''zbek tili'
$ws.Range("D42").Value = 'Kunduzgi'
$ws.Range("E42").Value = 'AD5077856'
$ws.Range("F42").Value = '50709066740028'
$ws.Range("G42").Value = 'Toshkent viloyati'
$ws.Range("H42").Value = 'Yuqori Chirchiq tumani'
$ws.Range("I42").Value = '998333515134'
$ws.Range("J42").Value = '+998770174575'
$ws.Range("K42").Value = '2025-05-14'

# Row 43
$ws.Range("A43").Value = 'Uaboyeva Sarvinoz Yusuf qizi'
$ws.Range("B43").Value = 'Yurisprudensiya'
$ws.Range("C43").Value = 'O''zbek tili'
$ws.Range("D43").Value = 'Kunduzgi'
$ws.Range("E43").Value = 'AD2635624'
$ws.Range("F43").Value = '63011066300059'
$ws.Range("G43").Value = 'Surxondaryo viloyati'
$ws.Range("H43").Value = 'Qumqoʻrgʻon tumani'
$ws.Range("I43").Value = '998941727977'
$ws.Range("J43").Value = '+998508890637'
$ws.Range("K43").Value = '2025-05-15'

# Row 44
$ws.Range("A44").Value = 'Xasan Saidmurodov'
$ws.Range("B44").Value = 'Yurisprudensiya'
$ws.Range("C44").Value = 'O''zbek tili'
$ws.Range("D44").Value = 'Kunduzgi'
$ws.Range("E44").Value = 'AD1999493'
$ws.Range("F44").Value = '52305066180114'
$ws.Range("G44").Value = 'Samarqand viloyati'
$ws.Range("H44").Value = 'Payariq tumani'
$ws.Range("I44").Value = '998957570177'
$ws.Range("J44").Value = '+998997787117'
$ws.Range("K44").Value = '2025-05-16'

# Row 45
$ws.Range("A45").Value = 'Abdurxmon Tuychibayev Abdudjabbar o''g''li'
$ws.Range("B45").Value = 'Bugalteriya hisobi'
$ws.Range("C45").Value = 'O''zbek tili'
$ws.Range("D45").Value = 'Kunduzgi'
$ws.Range("E45").Value = 'AD7006138'
$ws.Range("F45").Value = '51409076620044'
$ws.Range("G45").Value = 'Toshkent shahri'
$ws.Range("H45").Value = 'Yashnaobod tumani'
$ws.Range("I45").Value = '998882802529'
$ws.Range("J45").Value = '+998882802529'
$ws.Range("K45").Value = '2025-05-16'

# Row 46
$ws.Range("A46").Value = 'Rustamxojayev Javohirxoja muzaffarxon ogli'
$ws.Range("B46").Value = 'Yurisprudensiya'
$ws.Range("C46").Value = 'O''zbek tili'
$ws.Range("D46").Value = 'Kunduzgi'
$ws.Range("E46").Value = 'AD4246020'
$ws.Range("F46").Value = '52001076610011'
$ws.Range("G46").Value = 'Toshkent shahri'
$ws.Range("H46").Value = 'Olmazor tumani'
$ws.Range("I46").Value = '998974432526'
$ws.Range("J46").Value = '+998974432526'
$ws.Range("K46").Value = '2025-05-19'

# Row 47
$ws.Range("A47").Value = 'Nabiyev Sirojiddin Farohiddin o''g''li'
$ws.Range("B47").Value = 'Yurisprudensiya'
$ws.Range("C47").Value = 'O''zbek tili'
$ws.Range("D47").Value = 'Kunduzgi'
$ws.Range("E47").Value = 'AD6788338'
$ws.Range("F47").Value = '50102085150054'
$ws.Range("G47").Value = 'Andijon viloyati'
$ws.Range("H47").Value = 'Asaka tumani'
$ws.Range("I47").Value = '998507404744'
$ws.Range("J47").Value = '+998914854913'
$ws.Range("K47").Value = '2025-05-19'

# Row 48
$ws.Range("A48").Value = 'MAMASALIEV SOBIR BAXTIYOROVICH'
$ws.Range("B48").Value = 'Hayot faoliyati xavfsizligi'
$ws.Range("C48").Value = 'Rus tili'
$ws.Range("D48").Value = 'Kunduzgi'
$ws.Range("E48").Value = 'AC2787580'
$ws.Range("F48").Value = '52202047970013'
$ws.Range("G48").Value = 'Toshkent shahri'
$ws.Range("H48").Value = 'Yashnaobod tumani'
$ws.Range("I48").Value = '998944059922'
$ws.Range("J48").Value = '+998880132242'
$ws.Range("K48").Value = '2025-05-21'

# Restore default (non-text) formatting on the forced columns
$ws.Range("F35:F48").ClearFormats()
$ws.Range("I35:I48").ClearFormats()
$ws.Range("J35:J48").ClearFormats()
$ws.Range("K35:K48").ClearFormats()

